# Auto-generated Excel COM-interop script to apply odds updates
# described by the diff for Jogos_da_Semana_FlashScore_2024-11-20.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.38  # Odd_H_FT: 2.25 -> 2.38
$ws.Range("I2").Value = 3.25  # Odd_A_FT: 3.5 -> 3.25
$ws.Range("J2").Value = 3.25  # Odd_H_HT: 3.1 -> 3.25
$ws.Range("L2").Value = 4  # Odd_A_HT: 4.33 -> 4
$ws.Range("M2").Value = 1.11  # Odd_Over05_FT: 1.13 -> 1.11
$ws.Range("N2").Value = 6.5  # Odd_Under05_FT: 6 -> 6.5
$ws.Range("W2").Value = 6  # Odd_CS_1-0: 5.5 -> 6
$ws.Range("X2").Value = 10  # Odd_CS_2-0: 9.5 -> 10
$ws.Range("Z2").Value = 23  # Odd_CS_3-0: 21 -> 23
$ws.Range("AC2").Value = 6.5  # Odd_CS_0-0: 6 -> 6.5
$ws.Range("AE2").Value = 19  # Odd_CS_2-2: 21 -> 19
$ws.Range("AG2").Value = 7  # Odd_CS_0-1: 7.5 -> 7
$ws.Range("AJ2").Value = 34  # Odd_CS_0-3: 41 -> 34
$ws.Range("AL2").Value = 41  # Odd_CS_2-3: 51 -> 41
$ws.Range("AN2").Value = 4.33  # Odd_CS_1-0_HT: 4 -> 4.33
$ws.Range("AY2").Value = 34  # Odd_CS_1-2_HT: 41 -> 34
$ws.Range("AZ2").Value = 67  # Odd_CS_0-3_HT: 81 -> 67
$ws.Range("BB2").Value = 351  # Odd_CS_2-3_HT: 401 -> 351

# Row 4
$ws.Range("M4").Value = 1.06  # Odd_Over05_FT: 1.07 -> 1.06
$ws.Range("N4").Value = 10  # Odd_Under05_FT: 9 -> 10
$ws.Range("O4").Value = 1.33  # Odd_Over15_FT: 1.36 -> 1.33
$ws.Range("P4").Value = 3.25  # Odd_Under15_FT: 3 -> 3.25
$ws.Range("Q4").Value = 2.08  # Odd_Over25_FT: 2.1 -> 2.08
$ws.Range("R4").Value = 1.73  # Odd_Under25_FT: 1.7 -> 1.73
$ws.Range("BC4").Value = 151  # Odd_CS_3-3_HT: 126 -> 151
$ws.Range("BD4").Value = 151  # Odd_CS_4-4_HT: 126 -> 151

# Row 5
$ws.Range("O5").Value = 1.44  # Odd_Over15_FT: 1.5 -> 1.44
$ws.Range("P5").Value = 2.63  # Odd_Under15_FT: 2.5 -> 2.63
$ws.Range("Q5").Value = 2.4  # Odd_Over25_FT: 2.5 -> 2.4
$ws.Range("R5").Value = 1.53  # Odd_Under25_FT: 1.5 -> 1.53
$ws.Range("S5").Value = 1.53  # Odd_Over05_HT: 1.57 -> 1.53
$ws.Range("T5").Value = 2.38  # Odd_Under05_HT: 2.25 -> 2.38
$ws.Range("AT5").Value = 2.38  # Odd_CS_0-0_HT: 2.25 -> 2.38

# Row 8
$ws.Range("G8").Value = 2.15  # Odd_H_FT: 2.1 -> 2.15
$ws.Range("I8").Value = 3.3  # Odd_A_FT: 3.4 -> 3.3
$ws.Range("J8").Value = 3  # Odd_H_HT: 2.88 -> 3
$ws.Range("L8").Value = 4  # Odd_A_HT: 4.33 -> 4
$ws.Range("X8").Value = 9.5  # Odd_CS_2-0: 9 -> 9.5
$ws.Range("AE8").Value = 17  # Odd_CS_2-2: 19 -> 17
$ws.Range("AH8").Value = 15  # Odd_CS_0-2: 17 -> 15
$ws.Range("AI8").Value = 12  # Odd_CS_1-2: 13 -> 12
$ws.Range("AK8").Value = 29  # Odd_CS_1-3: 34 -> 29
$ws.Range("AO8").Value = 13  # Odd_CS_2-0_HT: 12 -> 13
$ws.Range("AW8").Value = 5  # Odd_CS_0-1_HT: 5.5 -> 5
$ws.Range("BB8").Value = 251  # Odd_CS_2-3_HT: 301 -> 251

# Row 9
$ws.Range("G9").Value = 1.65  # Odd_H_FT: 1.62 -> 1.65
$ws.Range("H9").Value = 3.7  # Odd_D_FT: 3.75 -> 3.7
$ws.Range("I9").Value = 5.5  # Odd_A_FT: 6 -> 5.5
$ws.Range("J9").Value = 2.25  # Odd_H_HT: 2.2 -> 2.25
$ws.Range("M9").Value = 1.05  # Odd_Over05_FT: 1.06 -> 1.05
$ws.Range("N9").Value = 11  # Odd_Under05_FT: 10 -> 11
$ws.Range("AI9").Value = 17  # Odd_CS_1-2: 19 -> 17
$ws.Range("AM9").Value = 251  # Odd_CS_4-4: 301 -> 251
$ws.Range("AN9").Value = 3.6  # Odd_CS_1-0_HT: 3.5 -> 3.6
$ws.Range("AO9").Value = 8.5  # Odd_CS_2-0_HT: 8 -> 8.5
$ws.Range("AQ9").Value = 29  # Odd_CS_3-0_HT: 26 -> 29
$ws.Range("AW9").Value = 6.5  # Odd_CS_0-1_HT: 7 -> 6.5

# Row 10
$ws.Range("I10").Value = 3.1  # Odd_A_FT: 3.2 -> 3.1
$ws.Range("J10").Value = 3.25  # Odd_H_HT: 3.2 -> 3.25
$ws.Range("M10").Value = 1.08  # Odd_Over05_FT: 1.1 -> 1.08
$ws.Range("N10").Value = 7.5  # Odd_Under05_FT: 7 -> 7.5
$ws.Range("AW10").Value = 4.75  # Odd_CS_0-1_HT: 5 -> 4.75
$ws.Range("AZ10").Value = 51  # Odd_CS_0-3_HT: 67 -> 51

# Row 11
$ws.Range("H11").Value = 3.3  # Odd_D_FT: 3.25 -> 3.3
$ws.Range("L11").Value = 4.33  # Odd_A_HT: 4 -> 4.33
$ws.Range("W11").Value = 6.5  # Odd_CS_1-0: 7 -> 6.5
$ws.Range("Y11").Value = 9  # Odd_CS_2-1: 9.5 -> 9
$ws.Range("Z11").Value = 19  # Odd_CS_3-0: 21 -> 19
$ws.Range("AG11").Value = 9.5  # Odd_CS_0-1: 9 -> 9.5
$ws.Range("AI11").Value = 13  # Odd_CS_1-2: 12 -> 13
$ws.Range("AW11").Value = 5.5  # Odd_CS_0-1_HT: 5 -> 5.5

# Row 12
$ws.Range("G12").Value = 3.1  # Odd_H_FT: 3.2 -> 3.1
$ws.Range("I12").Value = 2.25  # Odd_A_FT: 2.2 -> 2.25
$ws.Range("J12").Value = 3.5  # Odd_H_HT: 3.6 -> 3.5
$ws.Range("S12").Value = 1.33  # Odd_Over05_HT: 1.36 -> 1.33
$ws.Range("T12").Value = 3.25  # Odd_Under05_HT: 3 -> 3.25
$ws.Range("AD12").Value = 7  # Odd_CS_1-1: 6.5 -> 7
$ws.Range("AE12").Value = 12  # Odd_CS_2-2: 13 -> 12
$ws.Range("AG12").Value = 9.5  # Odd_CS_0-1: 9 -> 9.5
$ws.Range("AT12").Value = 3.25  # Odd_CS_0-0_HT: 3 -> 3.25
$ws.Range("AW12").Value = 4.5  # Odd_CS_0-1_HT: 4.33 -> 4.5

# Row 20
$ws.Range("G20").Value = 22  # Odd_H_FT: 25 -> 22
$ws.Range("H20").Value = 7.3  # Odd_D_FT: 7.4 -> 7.3
$ws.Range("I20").Value = 1.09  # Odd_A_FT: 1.08 -> 1.09
$ws.Range("J20").Value = 14  # Odd_H_HT: 16.5 -> 14
$ws.Range("K20").Value = 3.45  # Odd_D_HT: 3.55 -> 3.45
$ws.Range("L20").Value = 1.33  # Odd_A_HT: 1.3 -> 1.33
$ws.Range("S20").Value = 1.15  # Odd_Over05_HT: 1.14 -> 1.15
$ws.Range("T20").Value = 4.85  # Odd_Under05_HT: 5 -> 4.85
$ws.Range("U20").Value = 1.98  # Odd_BTTS_Yes: 2.02 -> 1.98
$ws.Range("V20").Value = 1.75  # Odd_BTTS_No: 1.7 -> 1.75
$ws.Range("W20").Value = 110  # Odd_CS_1-0: 120 -> 110
$ws.Range("X20").Value = 400  # Odd_CS_2-0: 500 -> 400
$ws.Range("Y20").Value = 90  # Odd_CS_2-1: 100 -> 90
$ws.Range("AA20").Value = 450  # Odd_CS_3-1: 600 -> 450
$ws.Range("AB20").Value = 200  # Odd_CS_3-2: 250 -> 200
$ws.Range("AC20").Value = 29  # Odd_CS_0-0: 28 -> 29
$ws.Range("AD20").Value = 20  # Odd_CS_1-1: 21 -> 20
$ws.Range("AE20").Value = 35  # Odd_CS_2-2: 37 -> 35
$ws.Range("AF20").Value = 110  # Odd_CS_3-3: 120 -> 110
$ws.Range("AG20").Value = 13.5  # Odd_CS_0-1: 13 -> 13.5
$ws.Range("AH20").Value = 8.25  # Odd_CS_0-2: 8 -> 8.25
$ws.Range("AI20").Value = 12  # Odd_CS_1-2: 12.5 -> 12
$ws.Range("AJ20").Value = 7.2  # Odd_CS_0-3: 7.1 -> 7.2
$ws.Range("AK20").Value = 10.5  # Odd_CS_1-3: 10.75 -> 10.5
$ws.Range("AL20").Value = 29  # Odd_CS_2-3: 30 -> 29
$ws.Range("AM20").Value = 600  # Odd_CS_4-4: 800 -> 600
$ws.Range("AN20").Value = 22  # Odd_CS_1-0_HT: 25 -> 22
$ws.Range("AO20").Value = 150  # Odd_CS_2-0_HT: 175 -> 150
$ws.Range("AP20").Value = 65  # Odd_CS_2-1_HT: 80 -> 65
$ws.Range("AT20").Value = 4.85  # Odd_CS_0-0_HT: 5 -> 4.85
$ws.Range("AU20").Value = 10  # Odd_CS_1-1_HT: 10.5 -> 10
$ws.Range("AV20").Value = 60  # Odd_CS_2-2_HT: 65 -> 60
$ws.Range("AW20").Value = 3.55  # Odd_CS_0-1_HT: 3.6 -> 3.55
$ws.Range("AX20").Value = 4.35  # Odd_CS_0-2_HT: 4.25 -> 4.35
$ws.Range("AZ20").Value = 7.6  # Odd_CS_0-3_HT: 7.2 -> 7.6
$ws.Range("BA20").Value = 21  # Odd_CS_1-3_HT: 20 -> 21
